# Append 29 new "Currency_Result" test-run rows (rows 55-83) recorded by a
# later run of the s2s Currency test suite. Each row is
# [TestData, Result, Status, TransactionID, Timestamp]; TransactionID is
# blank for the first few rows where the request failed before an id was
# issued.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Currency_Result")

$ws.Range("A55").Value = "INR"
$ws.Range("B55").Value = "FAIL"
$ws.Range("C55").Value = "Payment Failed"
$ws.Range("E55").Value = "2025-12-22 16:03:41"

$ws.Range("A56").Value = "INR"
$ws.Range("B56").Value = "FAIL"
$ws.Range("C56").Value = "Payment Failed"
$ws.Range("E56").Value = "2025-12-22 16:16:54"

$ws.Range("A57").Value = "INR"
$ws.Range("B57").Value = "FAIL"
$ws.Range("C57").Value = "Payment Failed"
$ws.Range("E57").Value = "2025-12-22 16:29:34"

$ws.Range("A58").Value = "INR"
$ws.Range("B58").Value = "FAIL"
$ws.Range("C58").Value = "Payment Failed"
$ws.Range("D58").Value = "69492674d8d85b68aa1abf16"
$ws.Range("E58").Value = "2025-12-22 16:37:48"

$ws.Range("A59").Value = "INR"
$ws.Range("B59").Value = "FAIL"
$ws.Range("C59").Value = "Payment Failed"
$ws.Range("D59").Value = "6949287ad8d85b68aa1ac73f"
$ws.Range("E59").Value = "2025-12-22 16:46:24"

$ws.Range("A60").Value = "INR"
$ws.Range("B60").Value = "FAIL"
$ws.Range("C60").Value = "Payment Failed"
$ws.Range("D60").Value = "69492e0c3fbd5815ff04ed0d"
$ws.Range("E60").Value = "2025-12-22 17:10:12"

$ws.Range("A61").Value = "INR"
$ws.Range("B61").Value = "FAIL"
$ws.Range("C61").Value = "Payment Failed"
$ws.Range("D61").Value = "69492f523fbd5815ff04f02e"
$ws.Range("E61").Value = "2025-12-22 17:15:34"

$ws.Range("A62").Value = "INR"
$ws.Range("B62").Value = "FAIL"
$ws.Range("C62").Value = "Payment Failed"
$ws.Range("D62").Value = "694930703fbd5815ff04f30e"
$ws.Range("E62").Value = "2025-12-22 17:20:20"

$ws.Range("A63").Value = "INR"
$ws.Range("B63").Value = "FAIL"
$ws.Range("C63").Value = "Payment Failed"
$ws.Range("D63").Value = "694931363fbd5815ff04f602"
$ws.Range("E63").Value = "2025-12-22 17:23:38"

$ws.Range("A64").Value = "INR"
$ws.Range("B64").Value = "FAIL"
$ws.Range("C64").Value = "Payment Failed"
$ws.Range("D64").Value = "694931853fbd5815ff04f8c3"
$ws.Range("E64").Value = "2025-12-22 17:26:39"

$ws.Range("A65").Value = "INR"
$ws.Range("B65").Value = "FAIL"
$ws.Range("C65").Value = "Payment Failed"
$ws.Range("D65").Value = "6949326d3fbd5815ff04fe5b"
$ws.Range("E65").Value = "2025-12-22 17:28:50"

$ws.Range("A66").Value = "USD"
$ws.Range("B66").Value = "FAIL"
$ws.Range("C66").Value = "Payment Failed"
$ws.Range("D66").Value = "694932e83fbd5815ff0501e9"
$ws.Range("E66").Value = "2025-12-22 17:30:53"

$ws.Range("A67").Value = "EUR"
$ws.Range("B67").Value = "FAIL"
$ws.Range("C67").Value = "Payment Failed"
$ws.Range("D67").Value = "694a58a0af0cbe08398a5fa4"
$ws.Range("E67").Value = "2025-12-23 14:24:07"

$ws.Range("A68").Value = "CLP"
$ws.Range("B68").Value = "FAIL"
$ws.Range("C68").Value = "Payment Failed"
$ws.Range("D68").Value = "694a78025f721d9eabe0c556"
$ws.Range("E68").Value = "2025-12-23 16:38:01"

$ws.Range("A69").Value = "COP"
$ws.Range("B69").Value = "FAIL"
$ws.Range("C69").Value = "Payment Failed"
$ws.Range("D69").Value = "694a7a215f721d9eabe0e23c"
$ws.Range("E69").Value = "2025-12-23 16:47:07"

$ws.Range("A70").Value = "CLP"
$ws.Range("B70").Value = "FAIL"
$ws.Range("C70").Value = "Payment Failed"
$ws.Range("D70").Value = "694a80076121018246f93d68"
$ws.Range("E70").Value = "2025-12-23 17:12:16"

$ws.Range("A71").Value = "CLP"
$ws.Range("B71").Value = "FAIL"
$ws.Range("C71").Value = "Payment Failed"
$ws.Range("D71").Value = "694a80c36121018246f93fec"
$ws.Range("E71").Value = "2025-12-23 17:15:21"

$ws.Range("A72").Value = "PEN"
$ws.Range("B72").Value = "FAIL"
$ws.Range("C72").Value = "Payment Failed"
$ws.Range("D72").Value = "694a860e6ab547d8d6261e4e"
$ws.Range("E72").Value = "2025-12-23 17:37:56"

$ws.Range("A73").Value = "PEN"
$ws.Range("B73").Value = "FAIL"
$ws.Range("C73").Value = "Payment Failed"
$ws.Range("D73").Value = "694a87136ab547d8d62620db"
$ws.Range("E73").Value = "2025-12-23 17:42:18"

$ws.Range("A74").Value = "PEN"
$ws.Range("B74").Value = "FAIL"
$ws.Range("C74").Value = "Payment Failed"
$ws.Range("D74").Value = "694a8a9ec89ddcf8f3998bec"
$ws.Range("E74").Value = "2025-12-23 17:57:29"

$ws.Range("A75").Value = "PEN"
$ws.Range("B75").Value = "FAIL"
$ws.Range("C75").Value = "Payment Failed"
$ws.Range("D75").Value = "694b7ce404f0708ed1e91269"
$ws.Range("E75").Value = "2025-12-24 11:11:11"

$ws.Range("A76").Value = "PEN"
$ws.Range("B76").Value = "FAIL"
$ws.Range("C76").Value = "Payment Failed"
$ws.Range("D76").Value = "694b7d5904f0708ed1e914f1"
$ws.Range("E76").Value = "2025-12-24 11:13:07"

$ws.Range("A77").Value = "PEN"
$ws.Range("B77").Value = "FAIL"
$ws.Range("C77").Value = "Payment Failed"
$ws.Range("D77").Value = "694b812f04f0708ed1e91c62"
$ws.Range("E77").Value = "2025-12-24 11:29:32"

$ws.Range("A78").Value = "PEN"
$ws.Range("B78").Value = "FAIL"
$ws.Range("C78").Value = "Payment Failed"
$ws.Range("D78").Value = "694b813204f0708ed1e91c68"
$ws.Range("E78").Value = "2025-12-24 11:29:38"

$ws.Range("A79").Value = "PEN"
$ws.Range("B79").Value = "FAIL"
$ws.Range("C79").Value = "Payment Failed"
$ws.Range("D79").Value = "694b841204f0708ed1e921fe"
$ws.Range("E79").Value = "2025-12-24 11:43:05"

$ws.Range("A80").Value = "CLP"
$ws.Range("B80").Value = "FAIL"
$ws.Range("C80").Value = "Payment Failed"
$ws.Range("D80").Value = "694b846904f0708ed1e92489"
$ws.Range("E80").Value = "2025-12-24 11:45:47"

$ws.Range("A81").Value = "MXN"
$ws.Range("B81").Value = "FAIL"
$ws.Range("C81").Value = "Payment Failed"
$ws.Range("D81").Value = "694b85e904f0708ed1e92a62"
$ws.Range("E81").Value = "2025-12-24 11:49:44"

$ws.Range("A82").Value = "INR"
$ws.Range("B82").Value = "FAIL"
$ws.Range("C82").Value = "Payment Failed"
$ws.Range("D82").Value = "694b91cd04f0708ed1e94191"
$ws.Range("E82").Value = "2025-12-24 12:40:21"

$ws.Range("A83").Value = "INR"
$ws.Range("B83").Value = "FAIL"
$ws.Range("C83").Value = "Payment Failed"
$ws.Range("D83").Value = "694b931804f0708ed1e94ce2"
$ws.Range("E83").Value = "2025-12-24 12:45:55"

